# The sheet receives a new weekly data row. It is inserted right after the
# previous week's row (between the current row 187 and row 188), which
# pushes every following row (old 188..261) down by one position (new
# 189..262) and extends the used range from A1:R261 to A1:R262.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 188; Excel shifts rows 188-261 down to
# 189-262 and grows the sheet dimension automatically.
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new weekly record.
$ws.Cells.Item(188, 1).Value = 9
$ws.Cells.Item(188, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(188, 3).Value = "Metropolitana"
$ws.Cells.Item(188, 4).Value = 44636
$ws.Cells.Item(188, 5).Value = 13
$ws.Cells.Item(188, 6).Value = 100112030
$ws.Cells.Item(188, 7).Value = "Poroto granado"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 52
$ws.Cells.Item(188, 11).Value = 22000
$ws.Cells.Item(188, 12).Value = 24000
$ws.Cells.Item(188, 13).Value = 23000
$ws.Cells.Item(188, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(188, 15).Value = "Provincia de Cardenal Caro"
$ws.Cells.Item(188, 16).Value = 920
$ws.Cells.Item(188, 17).Value = 25
$ws.Cells.Item(188, 18).Value = "Hortaliza"
